$d = $word.ActiveDocument

# --- Step 1: "Co the dung Command Pattern" (paragraph 7) -> two runs ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "FormHome: đổi màn hình(đổi chiến lược) khi chạy  và template method làm lại hoạt động "
$pEnd = $p7.Range.End
$rAfter = $d.Range($pEnd - 1, $pEnd - 1)
$rAfter.InsertAfter("(các child là các control)")
$newEnd = $p7.Range.End
$rSecond = $d.Range($pEnd - 1, $newEnd - 1)
# force a genuine run split even though formatting ends up identical
$rSecond.Bold = 1
$rSecond.Bold = 0

# --- Step 2: insert new red paragraph "Template method, strategy pattern" after paragraph 6 ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(7)
$pNew1.Range.Text = "Template method, strategy pattern"

# --- Step 3: "Singleton Pattern:" (paragraph 3) -> "MVC" ---
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "MVC"

# --- Step 4: insert new red paragraph "Singleton Pattern:" after paragraph 3 ---
$p3.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item(4)
$pNew2.Range.Text = "Singleton Pattern:"
